# Remove the trailing "Ver no Jupiter..." line, the copyright/footer line
# that follows it, and the blank paragraph that separates them from the
# preceding "Requisitos" text, per the site rebuild that dropped the
# Jekyll-site navigation/footer boilerplate from the exported document.

$d = $word.ActiveDocument

$verText  = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyText = "Contact: luizeleno@usp.br"

$count = $d.Paragraphs.Count
$verIdx  = -1
$copyIdx = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*$verText*")  { $verIdx  = $i }
    if ($t -like "*$copyText*") { $copyIdx = $i }
}

if ($verIdx -gt 0 -and $copyIdx -ge $verIdx) {
    # Also drop the blank paragraph immediately before the "Ver no Jupiter"
    # line, if there is one, so no stray empty paragraph is left behind.
    $startIdx = $verIdx
    if ($verIdx -gt 1) {
        $prevText = $d.Paragraphs.Item($verIdx - 1).Range.Text.Trim()
        if ($prevText -eq "") {
            $startIdx = $verIdx - 1
        }
    }

    $delStart = $d.Paragraphs.Item($startIdx).Range.Start
    $delEnd   = $d.Paragraphs.Item($copyIdx).Range.End

    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
